$wb = $excel.ActiveWorkbook

# --- Rework locale sheet names into proper .NET culture codes ---
$wsEn = $wb.Worksheets.Item("en")
$wsEn.Name = "en-US"
$wsRu = $wb.Worksheets.Item("ru")
$wsRu.Name = "ru-RU"

# --- Update the "main" sheet ---
$wsMain = $wb.Worksheets.Item("main")

# locales row now references the renamed locale codes
$wsMain.Range("B3").Value = "en-US"
$wsMain.Range("C3").Value = "ru-RU"

# Add a help note about culture codes, merged across D1:S1
$wsMain.Range("D1:S1").Merge()
$wsMain.Range("D1").Value = "You can find culture codes in https://lonewolfonline.net/list-net-culture-country-codes/"
$wsMain.Range("D1:S1").HorizontalAlignment = -4108
$wsMain.Range("D1:S1").VerticalAlignment = -4108

$wsMain.Range("C3").Select()

# --- Update selections on the other sheets to match the new workflow ---
$wsKeys = $wb.Worksheets.Item("keys")
$wsKeys.Range("B6").Select()

$wsSample = $wb.Worksheets.Item("sample_translation_list")
$wsSample.Range("C1").Select()

$wsRu.Range("B7").Select()

$wsEn.Range("B10").Select()
$wsEn.Activate()
